$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.093.33"
$ws.Range("E2").Value = '  +0.81%  '

$ws.Range("D3").Value = "'1.889.45"
$ws.Range("E3").Value = '  +0.13%  '

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").Value = "'0.7394"
$ws.Range("E5").Value = '  -0.92%  '

$ws.Range("D6").Value = "'242.59"
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = "'0.3168"
$ws.Range("E8").Value = '  +1.72%  '

$ws.Range("D9").Value = "'0.07203"
$ws.Range("E9").Value = '  +1.20%  '

$ws.Range("D10").Value = "'24.90"
$ws.Range("E10").Value = '  -1.55%  '

$ws.Range("E11").Value = '  -1.62%  '

$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = "'0.7595"
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'1.934.25"
$ws.Range("E13").Value = '  +2.15%  '

$ws.Range("D14").Value = "'5.438"
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").Value = "'92.77"
$ws.Range("E15").Value = '  -0.58%  '

$ws.Range("D16").Value = "'6.183"
$ws.Range("E16").Value = '  +0.62%  '

$ws.Range("D17").Value = "'30.170.92"
$ws.Range("E17").Value = '  +0.81%  '

$ws.Range("D18").Value = "'250.34"
$ws.Range("E18").Value = '  +2.92%  '

$ws.Range("D19").Value = "'13.64"
$ws.Range("E19").Value = '  -0.36%  '

$ws.Range("D20").Value = "'0.000007867"
$ws.Range("E20").Value = '  +0.93%  '

$ws.Range("D21").Value = "'2.169.81"
$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").Value = "'7.999"
$ws.Range("E23").Value = '  +0.13%  '

$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").Value = "'0.1582"
$ws.Range("E25").Value = '  -0.74%  '

$ws.Range("D26").Value = "'9.299"
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("D27").Value = "'164.22"
$ws.Range("E27").Value = '  +0.99%  '

$ws.Range("D28").Value = "'18.75"
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").Value = "'2.065"
$ws.Range("E29").Value = '  +1.92%  '

$ws.Range("D30").Value = "'1.480"
$ws.Range("E30").Value = '  -1.37%  '

$ws.Range("D31").Value = "'4.591"
$ws.Range("E31").Value = '  +2.62%  '

$ws.Range("E32").Value = '  +0.44%  '

$ws.Range("D33").Value = "'4.219"
$ws.Range("E33").Value = '  +2.85%  '

$ws.Range("D34").Value = "'0.05400"
$ws.Range("E34").Value = '  +0.26%  '

$ws.Range("D35").Value = "'1.252"
$ws.Range("E35").Value = '  +1.23%  '

$ws.Range("D36").Value = "'0.7704"
$ws.Range("E36").Value = '  +3.55%  '

$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("D38").Value = "'2.725"
$ws.Range("E38").Value = '  +0.59%  '

$ws.Range("D39").Value = "'0.01970"
$ws.Range("E39").Value = '  +1.97%  '

$ws.Range("D40").Value = "'2.764"
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").Value = "'0.4556"
$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("D42").Value = "'1.104.33"
$ws.Range("E42").Value = '  +1.21%  '

$ws.Range("D43").Value = "'6.080"
$ws.Range("E43").Value = '  +0.42%  '

$ws.Range("D44").Value = "'72.62"
$ws.Range("E44").Value = '  +0.16%  '

$ws.Range("D45").Value = "'0.8700"
$ws.Range("E45").Value = '  +1.60%  '

$ws.Range("D46").Value = "'104.46"
$ws.Range("E46").Value = '  +2.12%  '

$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("D48").Value = "'1.867"
$ws.Range("E48").Value = '  +0.41%  '

$ws.Range("D49").Value = "'7.600"
$ws.Range("E49").Value = '  -0.95%  '

$ws.Range("D50").Value = "'9.644"
$ws.Range("E50").Value = '  -0.47%  '

$ws.Range("D51").Value = "'2.073.25"
$ws.Range("E51").Value = '  +0.84%  '
